$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Some "Price" values are plain decimal numbers (e.g. "1.00", "9.50",
# "12.40"). Setting such a string straight into .Value makes Excel's COM
# automation auto-detect it as a number and normalize it (losing the
# trailing zero / exact text representation). To preserve the original
# text exactly - matching the source data, which stores these as plain
# strings - we force the cell's NumberFormat to Text ("@") immediately
# before assigning the value.

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.404.07"
$ws.Range("E2").Value = "  -2.03%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.487.98"
$ws.Range("E3").Value = "  -2.39%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.20"
$ws.Range("E5").Value = "  +4.61%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.05"
$ws.Range("E6").Value = "  -0.29%  "

# Row 7 - XRP (Price unchanged)
$ws.Range("E7").Value = "  -0.60%  "

# Row 8 - USDC
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.213"
$ws.Range("E9").Value = "  -4.74%  "

# Row 10 - Cardano
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.646"
$ws.Range("E10").Value = "  -0.96%  "

# Row 11 - Avalanche
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.94"
$ws.Range("E11").Value = "  -3.03%  "

# Row 12 - ShibaInu (Price unchanged)
$ws.Range("E12").Value = "  -3.89%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.50"
$ws.Range("E13").Value = "  +0.27%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.041.68"
$ws.Range("E14").Value = "  -2.41%  "

# Row 15 - BitcoinCash
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "599.78"
$ws.Range("E15").Value = "  +4.41%  "

# Row 16 - WrappedBTC
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.482.09"
$ws.Range("E16").Value = "  -1.94%  "

# Row 17 - Chainlink
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.82"
$ws.Range("E17").Value = "  -2.16%  "

# Row 18 - Uniswap
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.57"
$ws.Range("E18").Value = "  -1.87%  "

# Row 19 - WrappedEther
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.480.38"
$ws.Range("E19").Value = "  -2.34%  "

# Row 20 - TRON (Price unchanged)
$ws.Range("E20").Value = "  -0.20%  "

# Row 21 - Polygon
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.983"
$ws.Range("E21").Value = "  -1.95%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.16"
$ws.Range("E22").Value = "  -2.52%  "

# Row 23 - Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "106.16"
$ws.Range("E23").Value = "  +12.90%  "

# Row 24 - Toncoin (Price unchanged)
$ws.Range("E24").Value = "  +4.10%  "

# Row 25 - PancakeSwap (Price unchanged)
$ws.Range("E25").Value = "  +1.60%  "

# Row 26 - ImmutableX
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.04"
$ws.Range("E26").Value = "  +3.12%  "

# Row 27 - RenderToken
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.93"
$ws.Range("E27").Value = "  -2.68%  "

# Row 28 - Filecoin (Price unchanged)
$ws.Range("E28").Value = "  +4.90%  "

# Row 29 - EthereumClassic
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.39"
$ws.Range("E29").Value = "  +3.00%  "

# Row 30 - NEARProtocol
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.95"
$ws.Range("E30").Value = "  -3.37%  "

# Row 31 - dogwifhat (Price unchanged)
$ws.Range("E31").Value = "  +15.08%  "

# Row 32 - Cosmos
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.40"
$ws.Range("E32").Value = "  +0.81%  "

# Row 33 - Hedera (Price unchanged)
$ws.Range("E33").Value = "  -1.23%  "

# Row 34 - OKB
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.30"
$ws.Range("E34").Value = "  +0.30%  "

# Row 35 - Fetch.AI
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.17"
$ws.Range("E35").Value = "  -6.97%  "

# Row 36 - Dai (Price unchanged)
$ws.Range("E36").Value = "  -0.12%  "

# Row 37 - was Bittensor, now Stacks
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.69"
$ws.Range("E37").Value = "  +7.78%  "

# Row 38 - was Stacks, now Bittensor
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "520.21"
$ws.Range("E38").Value = "  -5.42%  "

# Row 39 - TheGraph
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.395"
$ws.Range("E39").Value = "  -4.63%  "

# Row 40 - Maker
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.587.96"
$ws.Range("E40").Value = "  +0.11%  "

# Row 41 - InjectiveProtocol
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.73"
$ws.Range("E41").Value = "  -3.67%  "

# Row 42 - PEPE
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0777"
$ws.Range("E42").Value = "  -3.37%  "

# Row 43 - Kaspa (Price unchanged)
$ws.Range("E43").Value = "  -0.93%  "

# Row 44 - VeChain
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0462"
$ws.Range("E44").Value = "  -0.31%  "

# Row 45 - ThetaToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.93"
$ws.Range("E45").Value = "  +0.20%  "

# Row 46 - Stellar (Price unchanged)
$ws.Range("E46").Value = "  +3.11%  "

# Row 47 - ApeXProtocol (Price unchanged)
$ws.Range("E47").Value = "  -4.69%  "

# Row 48 - THORChain
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.77"
$ws.Range("E48").Value = "  -6.22%  "

# Row 49 - FirstDigitalUSD (Price unchanged)
$ws.Range("E49").Value = "  +0.44%  "

# Row 50 - OceanProtocol (Price unchanged)
$ws.Range("E50").Value = "  -9.46%  "

# Row 51 - Monero
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.66"
$ws.Range("E51").Value = "  -2.78%  "
